# CobaltUsers.xlsx edit script
# - Rename Sheet2 -> Emails, populate it with two extra test accounts
# - Add 28 new user rows (53-80) to the Users sheet, each with a mailto hyperlink
# - Add 4 trailing spacer rows (81-84) with vertical-divider style in column E
# - Trim Sheet3 back down to a single data row
# - Re-point the active selection the way the author last left it

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Users sheet: new accounts
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Users")

$rowsData = @(
    @{ Row=53; User="SearchOpenWebUser1"; Email="SearchOpenWeb@mailinator.com " }
    @{ Row=54; User="FFHUser1"; Email="FFHUser1@mailinator.com " }
    @{ Row=55; User="FFHUser2"; Email="FFHUser2@mailinator.com" }
    @{ Row=56; User="FFHUser3"; Email="FFHUser3@mailinator.com" }
    @{ Row=57; User="FFHUser4"; Email="FFHUser4@mailinator.com" }
    @{ Row=58; User="FrontEndUser1"; Email="FrontEndUser1@mailinator.com" }
    @{ Row=59; User="FrontEndUser2"; Email="FrontEndUser2@mailinator.com" }
    @{ Row=60; User="FrontEndUser3"; Email="FrontEndUser3@mailinator.com" }
    @{ Row=61; User="FrontEndUser4"; Email="FrontEndUser4@mailinator.com" }
    @{ Row=62; User="FrontEndUser5"; Email="FrontEndUser5@mailinator.com" }
    @{ Row=63; User="FrontEndUser6"; Email="FrontEndUser6@mailinator.com" }
    @{ Row=64; User="FrontEndUser7"; Email="FrontEndUser7@mailinator.com" }
    @{ Row=65; User="FrontEndUser8"; Email="FrontEndUser8@mailinator.com" }
    @{ Row=66; User="FrontEndUser9"; Email="FrontEndUser9@mailinator.com" }
    @{ Row=67; User="FrontEndUser10"; Email="FrontEndUser10@mailinator.com" }
    @{ Row=68; User="UrlUser1"; Email="UrlUser1@mailinator.com" }
    @{ Row=69; User="UrlUser2"; Email="UrlUser2@mailinator.com" }
    @{ Row=70; User="UrlUser3"; Email="UrlUser3@mailinator.com" }
    @{ Row=71; User="LinkingUser1"; Email="LinkingUser1@mailinator.com " }
    @{ Row=72; User="LoginUser1"; Email="LoginUser1@mailinator.com " }
    @{ Row=73; User="LoginUser2"; Email="LoginUser2@mailinator.com " }
    @{ Row=74; User="LoginUser3"; Email="LoginUser3@mailinator.com " }
    @{ Row=75; User="LoginUser4"; Email="LoginUser4@mailinator.com " }
    @{ Row=76; User="LoginUser5"; Email="LoginUser5@mailinator.com " }
    @{ Row=77; User="LoginUser6"; Email="LoginUser6@mailinator.com " }
    @{ Row=78; User="LoginUser7"; Email="LoginUser7@mailinator.com " }
    @{ Row=79; User="CpetUser1"; Email="CpetUser1@mailinator.com " }
    @{ Row=80; User="CpetUser2"; Email="CpetUser2@mailinator.com " }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.User
    $ws.Cells.Item($row, 2).Value = "Password1"

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(7).Weight = 2
    $eCell.Borders.Item(10).LineStyle = 1
    $eCell.Borders.Item(10).Weight = 2

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = "N"
    $fCell.Borders.Item(7).LineStyle = 1
    $fCell.Borders.Item(7).Weight = 2
    $fCell.Borders.Item(10).LineStyle = 1
    $fCell.Borders.Item(10).Weight = 2

    $ws.Cells.Item($row, 7).Value = $r.Email
}

# Hyperlinks were historically added in this (non-sequential) order - replay it
# so the relationship ids line up the same way.
$hyperlinkOrder = @(55,54,56,57,58,59,60,61,62,63,64,65,66,67,68,69,53,70,71,72,73,74,75,76,77,78,79,80)
foreach ($row in $hyperlinkOrder) {
    $match = $rowsData | Where-Object { $_.Row -eq $row }
    $addr = "mailto:" + $match.Email.Trim()
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $addr)
}

# Trailing spacer rows, each with just a styled (bordered) blank E cell
for ($row = 81; $row -le 84; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(7).Weight = 2
    $eCell.Borders.Item(10).LineStyle = 1
    $eCell.Borders.Item(10).Weight = 2
}

# Column width tweaks (author widened columns A and E to fit the new content)
$ws.Columns.Item(1).ColumnWidth = 21.140625
$ws.Columns.Item(5).ColumnWidth = 39.28515625

# Restore the view: no frozen/scrolled topLeftCell, selection on C80
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C80").Select()

# ---------------------------------------------------------------------------
# 2. Rename Sheet2 -> Emails and populate it
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Emails"

$ws2.Range("A1").Value = "Email"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "tr-anz-tester1@yandex.com"
$ws2.Range("B2").Value = "tranztest"
$ws2.Range("A3").Value = "tr-anz-tester2@yandex.com"
$ws2.Range("B3").Value = "tranztest"

$ws2.Columns.Item(1).ColumnWidth = 26.28515625
$ws2.Columns.Item(2).ColumnWidth = 13.85546875

$ws2.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# 3. Sheet3: trim back down to a single row of data
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A3:A9").EntireRow.Delete()
